$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: A1 text changes from "Hola a todos" to "Hola" ---
$ws.Range("A1").Value = "Hola"

# --- Row 2: fill in the order that reproduces the shared-string table order
#     seen in the target workbook (D, I, J, A, C, E, F, G, B, H) ---
$ws.Range("D2").Value = "INDEFINIDO"
$ws.Range("I2").Value = "SIN_ASIGNAR"
$ws.Range("J2").Value = "No"

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "11-5-2025"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.0"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.0"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0.0"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.0"

$ws.Range("B2").Value = "18:47:19"

$ws.Range("H2").Value = "'"

# --- Rows 3-11: same pattern, all sharing B = "18:52:29" ---
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = "11-5-2025"
    $ws.Range("B$r").Value = "18:52:29"
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = "0.0"
    $ws.Range("D$r").Value = "INDEFINIDO"
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = "0.0"
    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = "0.0"
    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = "0.0"
    $ws.Range("H$r").Value = "'"
    $ws.Range("I$r").Value = "SIN_ASIGNAR"
    $ws.Range("J$r").Value = "No"
}

# --- Cursor / selection ends on B6, matching the saved workbook state ---
$ws.Range("B6").Select()
